$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert 5 fresh rows before the old "footer" row (row 13), pushing it
#    down to row 18. The newly inserted rows 13-17 come in blank/unstyled.
# ---------------------------------------------------------------------------
$ws.Rows("13:17").Insert()

# Give the new rows 13-17 the same look (borders/font) as the existing data
# rows (row 4 uses style s="6" for A:G and s="2" for H, which is what rows
# 13-17 need too).
$ws.Range("A4:H4").Copy()
$ws.Range("A13:H17").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 2. Update existing Code size / Data size values (C/D columns) for rows
#    3-8, 11-12, keeping the previous values around in new J/K "history"
#    columns for rows 3-6 and 12 (mirroring the style of each row).
# ---------------------------------------------------------------------------

# Row 3 (EventApp1): new code size 1500 (was 1426); keep old values in J3:K3.
$ws.Range("C3").Copy()
$ws.Range("J3:K3").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("J3").Value = 1426
$ws.Range("K3").Value = 1
$ws.Range("C3").Value = 1500

# Row 4 (EventApp2): new code size 978 (was 972); D4 unchanged (104).
$ws.Range("C4").Copy()
$ws.Range("J4:K4").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("J4").Value = 972
$ws.Range("K4").Value = 104
$ws.Range("C4").Value = 978

# Row 5 (EventApp3): new code size 1256 (was 1250); D5 unchanged (104).
$ws.Range("C5").Copy()
$ws.Range("J5:K5").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("J5").Value = 1250
$ws.Range("K5").Value = 104
$ws.Range("C5").Value = 1256

# Row 6 (EventApp4): new code size 1262 (was 1256); D6 unchanged (104).
$ws.Range("C6").Copy()
$ws.Range("J6:K6").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("J6").Value = 1256
$ws.Range("K6").Value = 104
$ws.Range("C6").Value = 1262

# Row 7 (FastPin1): code size 416 -> 380, data size 11 -> 7. No history cols.
$ws.Range("C7").Value = 380
$ws.Range("D7").Value = 7

# Row 8 (FastPin2): code size 180 -> 174. D8 unchanged (0).
$ws.Range("C8").Value = 174

# Row 11 (FastPin5): code size 362 -> 346. D11 unchanged (0).
$ws.Range("C11").Value = 346

# Row 12 (UartApp1): new code size 1884 (was 1618), new data size 169 (was 166).
$ws.Range("C12").Copy()
$ws.Range("J12:K12").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("J12").Value = 1618
$ws.Range("K12").Value = 166
$ws.Range("C12").Value = 1884
$ws.Range("D12").Value = 169

# ---------------------------------------------------------------------------
# 3. Fill in the 5 new example rows (13-17). The text values are entered in
#    this particular interleaved order so the shared-strings table ends up
#    built in the same sequence as the target workbook (UartApp2, SW UART,
#    UartApp3, SW UARX/UATX, PinChangeInterrupt1-3).
# ---------------------------------------------------------------------------
$ws.Range("A13").Value = "UartApp2"
$ws.Range("B14").Value = "SW UART"
$ws.Range("A14").Value = "UartApp3"
$ws.Range("B13").Value = "SW UARX/UATX"
$ws.Range("A15").Value = "PinChangeInterrupt1"
$ws.Range("A16").Value = "PinChangeInterrupt2"
$ws.Range("A17").Value = "PinChangeInterrupt3"

$ws.Range("C13").Value = 1628
$ws.Range("D13").Value = 158

$ws.Range("C14").Value = 1626
$ws.Range("D14").Value = 158

$ws.Range("C15").Value = 514
$ws.Range("D15").Value = 8

$ws.Range("C16").Value = 704
$ws.Range("D16").Value = 8

$ws.Range("C17").Value = 486
$ws.Range("D17").Value = 8

# ---------------------------------------------------------------------------
# 4. Misc sheet-level changes: selection, page orientation.
# ---------------------------------------------------------------------------
[void]$ws.Range("A18").Select()
$ws.PageSetup.Orientation = 2
